# Update mapa interactivo PEBCOM
# Remove the obsolete "4768 / VALLESE, FELIPE 684" record (row 23) and let
# every row below it shift up by one, reducing the used range from
# A1:P86 to A1:P85.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

$ws.Rows.Item(23).Delete()
